$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("E2").Value = 0.04

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("J3").Value = 0.15

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("J4").Value = 0.15

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 0.15
$ws.Range("K5").Value = 1

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("H6").Value = 0.24
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 0.26
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 2.02
